$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.3821797679500634
$ws.Range("D2").Value = 0.1164598247785733
$ws.Range("E2").Value = 0.1528805847024266
$ws.Range("F2").Value = 2.24439011372813
$ws.Range("G2").Value = 0.002527910193161421
$ws.Range("J2").Value = 0.2221476803079945
$ws.Range("K2").Value = 1.54036593527286
$ws.Range("L2").Value = 0.1311291549212932
$ws.Range("M2").Value = 0.5356190831377674
$ws.Range("N2").Value = 1.773140863759188
$ws.Range("O2").Value = 6.032816548646593

$ws.Range("C3").Value = 0.379915756642859
$ws.Range("D3").Value = 0.1151106274258282
$ws.Range("E3").Value = 0.153099355475895
$ws.Range("F3").Value = 2.251739254823214
$ws.Range("G3").Value = 0.002531104959684601
$ws.Range("J3").Value = 0.2234671281972531
$ws.Range("K3").Value = 1.439829800276584
$ws.Range("L3").Value = 0.1316593918146367
$ws.Range("M3").Value = 0.5162963364559019
$ws.Range("N3").Value = 1.78744045044936
$ws.Range("O3").Value = 6.060364373047605

$ws.Range("C4").Value = 0.3786914858254136
$ws.Range("D4").Value = 0.1143167553953006
$ws.Range("E4").Value = 0.1532838631964406
$ws.Range("F4").Value = 2.257336318726892
$ws.Range("G4").Value = 0.002533172424661997
$ws.Range("J4").Value = 0.2243520027899848
$ws.Range("K4").Value = 1.378461397391078
$ws.Range("L4").Value = 0.1320108852503346
$ws.Range("M4").Value = 0.5046084039483389
$ws.Range("N4").Value = 1.796819800191223
$ws.Range("O4").Value = 6.080305283952413

$ws.Range("C5").Value = 0.378234370772546
$ws.Range("D5").Value = 0.1140019811559441
$ws.Range("E5").Value = 0.1533716971467012
$ws.Range("F5").Value = 2.259889997019386
$ws.Range("G5").Value = 0.002534041634524889
$ws.Range("J5").Value = 0.2247314075392239
$ws.Range("K5").Value = 1.353545694184305
$ws.Range("L5").Value = 0.1321606558915907
$ws.Range("M5").Value = 0.4998901714195156
$ws.Range("N5").Value = 1.800792824462988
$ws.Range("O5").Value = 6.089192359844503

$ws.Range("C6").Value = 0.3781609933450909
$ws.Range("D6").Value = 0.1139502421695937
$ws.Range("E6").Value = 0.1533870464631946
$ws.Range("F6").Value = 2.260330514044284
$ws.Range("G6").Value = 0.002534187581192449
$ws.Range("J6").Value = 0.2247955441520908
$ws.Range("K6").Value = 1.349414089611599
$ws.Range("L6").Value = 0.1321859202830176
$ws.Range("M6").Value = 0.4991094225768862
$ws.Range("N6").Value = 1.801461658462969
$ws.Range("O6").Value = 6.090714013417909

$ws.Range("C7").Value = 0.3786851517135403
$ws.Range("D7").Value = 0.1143124748062903
$ws.Range("E7").Value = 0.1532849965183907
$ws.Range("F7").Value = 2.257369653743559
$ws.Range("G7").Value = 0.00253318403882744
$ws.Range("J7").Value = 0.2243570433826534
$ws.Range("K7").Value = 1.378124999185957
$ws.Range("L7").Value = 0.1320128786326311
$ws.Range("M7").Value = 0.5045445907691644
$ws.Range("N7").Value = 1.796872770743711
$ws.Range("O7").Value = 6.080422056983679

$ws.Range("C8").Value = 0.3813647647090335
$ws.Range("D8").Value = 0.1159874817022626
$ws.Range("E8").Value = 0.1529456176333497
$ws.Range("F8").Value = 2.246699021549659
$ws.Range("G8").Value = 0.002528989824790048
$ws.Range("J8").Value = 0.2225871281235392
$ws.Range("K8").Value = 1.505627033281115
$ws.Range("L8").Value = 0.1313066086595445
$ws.Range("M8").Value = 0.5289202151068011
$ws.Range("N8").Value = 1.777947058070154
$ws.Range("O8").Value = 6.041686875597577

$ws.Range("C9").Value = 0.3879323826221253
$ws.Range("D9").Value = 0.1195442998297622
$ws.Range("E9").Value = 0.1526771733792565
$ws.Range("F9").Value = 2.234378449173136
$ws.Range("G9").Value = 0.0025216013044854
$ws.Range("J9").Value = 0.2197084029779397
$ws.Range("K9").Value = 1.7584672173939
$ws.Range("L9").Value = 0.1301266657612761
$ws.Range("M9").Value = 0.5781066983801608
$ws.Range("N9").Value = 1.745582574060997
$ws.Range("O9").Value = 5.989746281364177

$ws.Range("C10").Value = 0.3935550052284498
$ws.Range("D10").Value = 0.1223211653692005
$ws.Range("E10").Value = 0.1527206875282872
$ws.Range("F10").Value = 2.23057234958128
$ws.Range("G10").Value = 0.002516677657417078
$ws.Range("J10").Value = 0.217953223115174
$ws.Range("K10").Value = 1.94588462307803
$ws.Range("L10").Value = 0.1293838750310758
$ws.Range("M10").Value = 0.615075593061448
$ws.Range("N10").Value = 1.724690212224075
$ws.Range("O10").Value = 5.966242035779402

$ws.Range("C11").Value = 0.3962854641870308
$ws.Range("D11").Value = 0.1236195251715628
$ws.Range("E11").Value = 0.1527924874257671
$ws.Range("F11").Value = 2.229980229942072
$ws.Range("G11").Value = 0.002514546259199084
$ws.Range("J11").Value = 0.2172326467223193
$ws.Range("K11").Value = 2.031494405899991
$ws.Range("L11").Value = 0.1290727229733371
$ws.Range("M11").Value = 0.632071715858963
$ws.Range("N11").Value = 1.715810533304513
$ws.Range("O11").Value = 5.958735405497322

$ws.Range("C12").Value = 0.3973441849697679
$ws.Range("D12").Value = 0.1241161899659318
$ws.Range("E12").Value = 0.1528271299198387
$ws.Range("F12").Value = 2.22991982364627
$ws.Range("G12").Value = 0.002513754658181266
$ws.Range("J12").Value = 0.2169709617839288
$ws.Range("K12").Value = 2.063961980883221
$ws.Range("L12").Value = 0.128958729217894
$ws.Range("M12").Value = 0.6385331112120198
$ws.Range("N12").Value = 1.712537690602396
$ws.Range("O12").Value = 5.956351133628459

$ws.Range("C13").Value = 0.3971150707630784
$ws.Range("D13").Value = 0.1240090025566332
$ws.Range("E13").Value = 0.1528193379397287
$ws.Range("F13").Value = 2.229925547622116
$ws.Range("G13").Value = 0.002513924454796092
$ws.Range("J13").Value = 0.2170268232454546
$ws.Range("K13").Value = 2.056967355116853
$ws.Range("L13").Value = 0.1289831095620766
$ws.Range("M13").Value = 0.6371404134299823
$ws.Range("N13").Value = 1.713238567715599
$ws.Range("O13").Value = 5.95684424110928

$ws.Range("C14").Value = 0.3963720701018758
$ws.Range("D14").Value = 0.1236602860808063
$ws.Range("E14").Value = 0.1527951882486178
$ws.Range("F14").Value = 2.229971977219222
$ws.Range("G14").Value = 0.002514480823169686
$ws.Range("J14").Value = 0.2172108937322079
$ws.Range("K14").Value = 2.034164560380418
$ws.Range("L14").Value = 0.1290632679004062
$ws.Range("M14").Value = 0.632602793196277
$ws.Range("N14").Value = 1.715539477321357
$ws.Range("O14").Value = 5.958530063060067

$ws.Range("C15").Value = 0.3959201822718228
$ws.Range("D15").Value = 0.1234473373626201
$ws.Range("E15").Value = 0.1527813657871206
$ws.Range("F15").Value = 2.230021750089918
$ws.Range("G15").Value = 0.002514823633017924
$ws.Range("J15").Value = 0.2173250979662953
$ws.Range("K15").Value = 2.020203522010206
$ws.Range("L15").Value = 0.1291128659405061
$ws.Range("M15").Value = 0.6298266576600398
$ws.Range("N15").Value = 1.716960529972887
$ws.Range("O15").Value = 5.959622372613069

$ws.Range("C16").Value = 0.3933800317344662
$ws.Range("D16").Value = 0.1222370174955429
$ws.Range("E16").Value = 0.1527170395854434
$ws.Range("F16").Value = 2.230633970399396
$ws.Range("G16").Value = 0.002516819124277979
$ws.Range("J16").Value = 0.2180018800068133
$ws.Range("K16").Value = 1.940296773108628
$ws.Range("L16").Value = 0.1294047466206543
$ws.Range("M16").Value = 0.6139684218754695
$ws.Range("N16").Value = 1.725283076606075
$ws.Range("O16").Value = 5.966796734378505

$ws.Range("C17").Value = 0.3918659080752036
$ws.Range("D17").Value = 0.1215034922503477
$ws.Range("E17").Value = 0.1526908795832043
$ws.Range("F17").Value = 2.231301338943936
$ws.Range("G17").Value = 0.002518071003109742
$ws.Range("J17").Value = 0.2184369958962016
$ws.Range("K17").Value = 1.891365682098524
$ws.Range("L17").Value = 0.1295906469668253
$ws.Range("M17").Value = 0.6042854309991554
$ws.Range("N17").Value = 1.7305485420991
$ws.Range("O17").Value = 5.972014046759369

$ws.Range("C18").Value = 0.391011284207039
$ws.Range("D18").Value = 0.121084899391704
$ws.Range("E18").Value = 0.1526807294846932
$ws.Range("F18").Value = 2.231792428089094
$ws.Range("G18").Value = 0.002518801257519507
$ws.Range("J18").Value = 0.218694592834396
$ws.Range("K18").Value = 1.863255108740645
$ws.Range("L18").Value = 0.1297000902966587
$ws.Range("M18").Value = 0.5987328814934116
$ws.Range("N18").Value = 1.733635871022294
$ws.Range("O18").Value = 5.975314732714821

$ws.Range("C19").Value = 0.3907247183265667
$ws.Range("D19").Value = 0.1209437413266414
$ws.Range("E19").Value = 0.1526781345633381
$ws.Range("F19").Value = 2.231977120802995
$ws.Range("G19").Value = 0.002519050264535077
$ws.Range("J19").Value = 0.2187830700903035
$ws.Range("K19").Value = 1.853743122268952
$ws.Range("L19").Value = 0.1297375788705288
$ws.Range("M19").Value = 0.5968557892728441
$ws.Range("N19").Value = 1.734691284412534
$ws.Range("O19").Value = 5.976483776263251

$ws.Range("C20").Value = 0.3920254068446525
$ws.Range("D20").Value = 0.1215812349070831
$ws.Range("E20").Value = 0.1526931577771009
$ws.Range("F20").Value = 2.231219197948178
$ws.Range("G20").Value = 0.00251793668257839
$ws.Range("J20").Value = 0.2183899185707112
$ws.Range("K20").Value = 1.896571046620977
$ws.Range("L20").Value = 0.1295705970136574
$ws.Range("M20").Value = 0.6053144606534602
$ws.Range("N20").Value = 1.729981942201086
$ws.Range("O20").Value = 5.971427622320675

$ws.Range("C21").Value = 0.3965896361833643
$ws.Range("D21").Value = 0.1237625771532862
$ws.Range("E21").Value = 0.1528020795115985
$ws.Range("F21").Value = 2.229953893879951
$ws.Range("G21").Value = 0.00251431698387206
$ws.Range("J21").Value = 0.2171565244187121
$ws.Range("K21").Value = 2.040860977429759
$ws.Range("L21").Value = 0.1290396195394816
$ws.Range("M21").Value = 0.6339349177323115
$ws.Range("N21").Value = 1.71486121062609
$ws.Range("O21").Value = 5.958022455475117

$ws.Range("C22").Value = 0.3997168799570545
$ws.Range("D22").Value = 0.1252173549736142
$ws.Range("E22").Value = 0.15291669902523
$ws.Range("F22").Value = 2.230081787950752
$ws.Range("G22").Value = 0.002512041690742478
$ws.Range("J22").Value = 0.2164155987419534
$ws.Range("K22").Value = 2.135447658727799
$ws.Range("L22").Value = 0.1287149292521992
$ws.Range("M22").Value = 0.6527875242283727
$ws.Range("N22").Value = 1.705501738231249
$ws.Range("O22").Value = 5.951932995957691

$ws.Range("C23").Value = 0.3980346376298769
$ws.Range("D23").Value = 0.1244382625394564
$ws.Range("E23").Value = 0.1528515584974386
$ws.Range("F23").Value = 2.229926164042809
$ws.Range("G23").Value = 0.002513247810931386
$ws.Range("J23").Value = 0.216805086685504
$ws.Range("K23").Value = 2.084939462845398
$ws.Range("L23").Value = 0.1288861834037878
$ws.Range("M23").Value = 0.6427121627886976
$ws.Range("N23").Value = 1.710449254626639
$ws.Range("O23").Value = 5.954938522228503

$ws.Range("C24").Value = 0.3919532480250751
$ws.Range("D24").Value = 0.1215460777286168
$ws.Range("E24").Value = 0.1526921125731953
$ws.Range("F24").Value = 2.231255999310406
$ws.Range("G24").Value = 0.002517997376027662
$ws.Range("J24").Value = 0.2184111790440646
$ws.Range("K24").Value = 1.894217637138183
$ws.Range("L24").Value = 0.1295796536012528
$ws.Range("M24").Value = 0.6048491916783405
$ws.Range("N24").Value = 1.730237914662681
$ws.Range("O24").Value = 5.97169180664028

$ws.Range("C25").Value = 0.3860154223818739
$ws.Range("D25").Value = 0.1185531835547238
$ws.Range("E25").Value = 0.1527074226922238
$ws.Range("F25").Value = 2.236790295248724
$ws.Range("G25").Value = 0.002523511094935467
$ws.Range("J25").Value = 0.2204239058841502
$ws.Range("K25").Value = 1.689772408408828
$ws.Range("L25").Value = 0.1304240119785423
$ws.Range("M25").Value = 0.564653458386509
$ws.Range("N25").Value = 1.753830604800243
$ws.Range("O25").Value = 6.001224772901935
